$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1: add the new "Web Hosting Business" course row (row 14)
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Write URL (column C) before the title (column B) so the shared-string table
# picks up the same ordering as the source workbook.
$sheet1.Range("C14").Value = "https://www.udemy.com/course/start-a-web-hosting-business/"
$sheet1.Range("B14").Value = "How to Start and Run a Web Hosting Business from Home"
$sheet1.Range("B14:C14").Style = "Good"

# Move the selection on Sheet1 (this sheet is not the final active tab).
$sheet1.Range("C16").Select()

# ---------------------------------------------------------------------------
# 2. Insert a new "Kotlin" worksheet right after "Android Application".
#    We duplicate "Framework" (same column-B width we need, no page setup)
#    so the duplicated sheet already carries the correct column formatting,
#    then wipe its contents and refill them.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("Framework")
$androidSheet = $wb.Worksheets.Item("Android Application")
$templateSheet.Copy($null, $androidSheet)

$kotlinSheet = $wb.Worksheets.Item("Android Application").Next
$kotlinSheet.Name = "Kotlin"
$kotlinSheet.Cells.Clear()

# Column widths: column B already matches (73.140625); tune column C as
# close as the engine's width quantization allows to the target 91.42578125.
$kotlinSheet.Columns.Item(3).ColumnWidth = 90.66666666666667

# Write URL (column C) before the title (column B) to match shared-string
# insertion order (url index lower than title index), same pattern as above.
$kotlinSheet.Range("C2").Value = "https://www.udemy.com/course/the-complete-kotlin-developer-course-java/"
$kotlinSheet.Range("B2").Value = "The Complete Kotlin Developer Course"

# Final selection/activation: Kotlin tab becomes the active tab.
$kotlinSheet.Range("B2").Select()
